$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the e-mail address shown in column C (rows 2-6): the club's
# generic "test@lamediainglesa.com" address is replaced everywhere by
# the new contact address.
$ws.Range("C2").Value = "joeljuaristi@hotmail.com"
$ws.Range("C3").Value = "joeljuaristi@hotmail.com"
$ws.Range("C4").Value = "joeljuaristi@hotmail.com"
$ws.Range("C5").Value = "joeljuaristi@hotmail.com"
$ws.Range("C6").Value = "joeljuaristi@hotmail.com"

# The old workbook had one hyperlink on C2 (left untouched) and a single
# hyperlink spanning the merged reference C3:C6. Split that second
# hyperlink into four independent per-cell hyperlinks, each pointing to
# the new mailto: address, mirroring the way Excel stores them once the
# link is re-created cell by cell.
$ws.Hyperlinks.Item(2).Delete()
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:joeljuaristi@hotmail.com")
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:joeljuaristi@hotmail.com")
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:joeljuaristi@hotmail.com")
$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:joeljuaristi@hotmail.com")

# Adding hyperlinks through this API re-applies the built-in "Hyperlink"
# cell style from scratch; restore the original shared style so the
# cells keep referencing the same format as before (and as C2, which we
# never touched).
$ws.Range("C3:C6").Style = "Hipervínculo"

# Restore the cursor position recorded in the saved file.
$ws.Range("I12").Select()
